$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '44.267.94'
$ws.Range("E2").Value = '  +2.09%  '
$ws.Range("D3").Value = '2.374.37'
$ws.Range("E3").Value = '  -0.02%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.692'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +7.31%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '244.76'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.46%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '76.87'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +5.48%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.596'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +23.25%  '
$ws.Range("E10").Value = '  +6.14%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '57.91'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.84%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '32.23'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +18.64%  '
$ws.Range("E13").Value = '  +19.35%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.108'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.23%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '17.22'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +7.29%  '
$ws.Range("D16").Value = '2.726.80'
$ws.Range("E16").Value = '  +0.00%  '
$ws.Range("E17").Value = '  +7.26%  '
$ws.Range("D18").Value = '2.370.14'
$ws.Range("E18").Value = '  +0.03%  '
$ws.Range("D19").Value = '44.299.60'
$ws.Range("E19").Value = '  +2.13%  '
$ws.Range("E20").Value = '  +4.27%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.74'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +6.01%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '78.80'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +5.54%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '258.23'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.06%  '
$ws.Range("E24").Value = '  +0.05%  '
$ws.Range("E25").Value = '  +4.60%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.72'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.83%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.92'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +8.88%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.79'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +17.40%  '
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.30'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.69%  '
$ws.Range("B30").Value = 'EthereumClassic'
$ws.Range("C30").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '23.16'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.75%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '175.16'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.57%  '
$ws.Range("E32").Value = '  +1.26%  '
$ws.Range("E33").Value = '  +7.00%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.40'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +7.92%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0762'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +10.02%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.34'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.47%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.90'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +5.22%  '
$ws.Range("E38").Value = '  +1.23%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.60'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.52%  '
$ws.Range("E40").Value = '  +8.60%  '
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '9.17'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.83%  '
$ws.Range("B42").Value = 'InjectiveProtocol'
$ws.Range("C42").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '19.16'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.98%  '
$ws.Range("E43").Value = '  -0.10%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.194'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +16.17%  '
$ws.Range("E45").Value = '  +2.82%  '
$ws.Range("E46").Value = '  +4.85%  '
$ws.Range("E47").Value = '  +5.76%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.54'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +12.62%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '102.54'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.36%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.48'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.68%  '
$ws.Range("D51").Value = '1.476.76'
$ws.Range("E51").Value = '  +2.34%  '
